# Add Forgot Password and Reset Password
# Replace the placeholder demo user list (10 fake rows) with the real
# user list (3 rows: admin + 2 users), update their hyperlinked emails,
# and switch the workbook font from Calibri to Arial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all existing mailto hyperlinks before we touch the rows/cells they
# are anchored to, and remove the now-unused trailing rows (4-10) - only
# 3 real users remain.
$ws.Hyperlinks.Delete()
$ws.Rows("4:10").Delete()

# --- Row 1: Le Dinh Cuong -------------------------------------------------
$ws.Range("A1").Value = "Lê Đình Cường"
$ws.Range("B1").Value = "dinhcuong1.firewin99@gmail.com"

# --- Row 2: Vu Xuan Canh --------------------------------------------------
$ws.Range("A2").Value = "Vũ Xuân Cảnh"
$ws.Range("B2").Value = "xuancanhit99@gmail.com"

# --- Row 3: Phuong Tien Dong ----------------------------------------------
$ws.Range("A3").Value = "Phương Tiến Đông"
$ws.Range("B3").Value = "dongpt4101@gmail.com"

# Re-create the mailto: hyperlinks on the (now correct) e-mail cells.
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:dinhcuong1.firewin99@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:xuancanhit99@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:dongpt4101@gmail.com")

# Switch the workbook's font from Calibri to Arial everywhere.
$ws.Cells.Font.Name = "Arial"

# Restore the last-used selection.
$ws.Range("B5").Select()
